## Timesheet update: add three new log rows (7-9) for the "Programming" work
## on collision detection / hit reactions, and the tilemap / Monogame Extended
## troubleshooting session — per commit "created new solution to work out
## issues with Monogame Extended".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7 -------------------------------------------------------------
# Seed A7:D7 from an existing row so the date/time cells pick up the SAME
# style indices already in the workbook (s=1 date, s=2 time) instead of the
# engine minting brand-new (duplicate) cellXfs entries.
$ws.Range("A5:D5").Copy($ws.Range("A7:D7"))
$ws.Range("A7").Value = 43898
$ws.Range("B7").Value = 0.45833333333333331
$ws.Range("C7").Value = 0.66666666666666663
$ws.Range("D7").Value = 0.020833333333333332
$ws.Range("E7").Value = "4 hours 30 minutes"
$ws.Range("F7").Value = "Programming"
$ws.Range("G7").Value = "Added code for collision detection and created attack animation for player"

# --- Row 8 -------------------------------------------------------------
$ws.Range("A5:D5").Copy($ws.Range("A8:D8"))
$ws.Range("A8").Value = 43902
$ws.Range("B8").Value = 0.5
$ws.Range("C8").Value = 0.625
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = "3 hours"
$ws.Range("F8").Value = "Proogramming"
$ws.Range("G8").Value = "Added hit detection and reactions for the player's attacks"

# --- Row 9 -------------------------------------------------------------
$ws.Range("A5:D5").Copy($ws.Range("A9:D9"))
$ws.Range("A9").Value = 43905
$ws.Range("B9").Value = 0.45833333333333331
# This entry's start-time cell carries a distinct time format (h:mm:ss AM/PM)
# from the rest of the sheet -- giving the workbook its 4th cellXfs entry.
$ws.Range("B9").NumberFormat = "h:mm:ss AM/PM"
$ws.Range("C9").Value = 0.625
$ws.Range("D9").Value = 0.020833333333333332
$ws.Range("E9").Value = "3 hours 30 minutes"
$ws.Range("F9").Value = "Research/Programming"
$ws.Range("G9").Value = "Created tilemap for starting room and attempted unsuccessfully to get monogame extended library working"

# --- Cosmetic follow-ups ------------------------------------------------
# Column widths grew to fit the new (longer) activity/comments text.
$ws.Columns.Item(2).ColumnWidth = 10.592447916666666
$ws.Columns.Item(3).ColumnWidth = 11.877604166666666
$ws.Columns.Item(7).ColumnWidth = 98.307291666666671

# Selection left where the author last clicked when saving.
[void]$ws.Range("G14").Select()
